$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 writes
$ws.Cells.Item(2, 4).Value = 0.04155
$ws.Cells.Item(2, 5).Value = 0.0546
$ws.Cells.Item(2, 7).Value = 0.2128576072821846
$ws.Cells.Item(2, 8).Value = 0.2115572171651495
$ws.Cells.Item(2, 9).Value = 0.383777633289987
$ws.Cells.Item(2, 10).Value = 0.3564938545847903
$ws.Cells.Item(2, 11).Value = 5.835
$ws.Cells.Item(2, 12).Value = 0.4742360208062419
$ws.Cells.Item(2, 13).Value = 3.226
$ws.Cells.Item(2, 14).Value = 0.0763369616658779
$ws.Cells.Item(2, 15).Value = 0.5528706083976007
$ws.Cells.Item(2, 16).Value = 3.226
$ws.Cells.Item(2, 17).Value = 0.0763369616658779
$ws.Cells.Item(2, 18).Value = 0.5528706083976007
$ws.Cells.Item(2, 19).Value = 0.0
$ws.Cells.Item(2, 20).Value = 0.0
$ws.Cells.Item(2, 21).Value = 6.038
$ws.Cells.Item(2, 22).Value = 0.1428774254614293
$ws.Cells.Item(2, 23).Value = 0.09771689497716894
$ws.Cells.Item(2, 24).Value = 0.06987381421014414
$ws.Cells.Item(2, 25).Value = 0.0278430807670248
$ws.Cells.Item(2, 26).Value = 0.3138375207243974
$ws.Cells.Item(2, 28).Value = 0.06987381421014414
$ws.Cells.Item(2, 29).Value = -0.06987381421014414
$ws.Cells.Item(2, 30).Value = 5.832
$ws.Cells.Item(2, 32).Value = 5.832
$ws.Cells.Item(2, 33).Value = -0.2060000000000004
$ws.Cells.Item(2, 34).Value = 0.1212675704898944
$ws.Cells.Item(2, 35).Value = 0.1131877729257642
$ws.Cells.Item(2, 36).Value = -0.004898463879773635
$ws.Cells.Item(2, 37).Value = -0.004528766460747036
$ws.Cells.Item(2, 38).Value = 0.509
$ws.Cells.Item(2, 39).Value = 0.398
$ws.Cells.Item(2, 40).Value = 1.158982511923688
$ws.Cells.Item(2, 41).Value = 9.277013752455796
$ws.Cells.Item(2, 42).Value = -0.04093799682034984
$ws.Cells.Item(2, 43).Value = 11.8643216080402

# Row 3 writes
$ws.Cells.Item(3, 4).Value = 0.0313
$ws.Cells.Item(3, 5).Value = 0.0546
$ws.Cells.Item(3, 7).Value = 0.4129263913824057
$ws.Cells.Item(3, 8).Value = 0.4129263913824057
$ws.Cells.Item(3, 9).Value = 0.5080789946140036
$ws.Cells.Item(3, 10).Value = 0.3953146997662194
$ws.Cells.Item(3, 11).Value = 2.16
$ws.Cells.Item(3, 12).Value = 0.3877917414721724
$ws.Cells.Item(3, 13).Value = 1.02
$ws.Cells.Item(3, 14).Value = 0.08095238095238096
$ws.Cells.Item(3, 15).Value = 0.4722222222222222
$ws.Cells.Item(3, 16).Value = 1.02
$ws.Cells.Item(3, 17).Value = 0.08095238095238096
$ws.Cells.Item(3, 18).Value = 0.4722222222222222
$ws.Cells.Item(3, 21).Value = 0.16
$ws.Cells.Item(3, 22).Value = 0.0126984126984127
$ws.Cells.Item(3, 23).Value = 0.2288135593220339
$ws.Cells.Item(3, 24).Value = 0.07088001256841266
$ws.Cells.Item(3, 25).Value = 0.1579335467536213
$ws.Cells.Item(3, 26).Value = 0.6437817845584837
$ws.Cells.Item(3, 27).Value = 0.2544964028776979
$ws.Cells.Item(3, 28).Value = 0.07046270422986359
$ws.Cells.Item(3, 29).Value = 0.1840336986478343
$ws.Cells.Item(3, 30).Value = 0.322
$ws.Cells.Item(3, 32).Value = 0.322
$ws.Cells.Item(3, 33).Value = 0.162
$ws.Cells.Item(3, 34).Value = 0.02491874322860239
$ws.Cells.Item(3, 35).Value = 0.03119550474714203
$ws.Cells.Item(3, 36).Value = 0.01269393511988716
$ws.Cells.Item(3, 37).Value = 0.01594174375123007
$ws.Cells.Item(3, 38).Value = 0.045
$ws.Cells.Item(3, 39).Value = 0.045
$ws.Cells.Item(3, 40).Value = 0.1066225165562914
$ws.Cells.Item(3, 41).Value = 62.88888888888889
$ws.Cells.Item(3, 42).Value = 0.05364238410596027
$ws.Cells.Item(3, 43).Value = 62.88888888888889

# Row 4 writes
$ws.Cells.Item(4, 2).Value = "Cyan Limited (KASE:CYAN)"
$ws.Cells.Item(4, 4).Value = -0.139
$ws.Cells.Item(4, 5).Value = -0.189
$ws.Cells.Item(4, 7).Value = 0.2710144927536232
$ws.Cells.Item(4, 8).Value = 0.2652173913043478
$ws.Cells.Item(4, 9).Value = 0.6666666666666667
$ws.Cells.Item(4, 10).Value = 0.589029535864979
$ws.Cells.Item(4, 11).Value = 1.4
$ws.Cells.Item(4, 12).Value = 0.5072463768115942
$ws.Cells.Item(4, 13).Value = 0.361
$ws.Cells.Item(4, 14).Value = 0.02542253521126761
$ws.Cells.Item(4, 15).Value = 0.2578571428571428
$ws.Cells.Item(4, 16).Value = 0.361
$ws.Cells.Item(4, 17).Value = 0.02542253521126761
$ws.Cells.Item(4, 18).Value = 0.2578571428571428
$ws.Cells.Item(4, 21).Value = 0.029
$ws.Cells.Item(4, 22).Value = 0.002042253521126761
$ws.Cells.Item(4, 23).Value = 0.1217391304347826
$ws.Cells.Item(4, 24).Value = 0.07982796881995903
$ws.Cells.Item(4, 25).Value = 0.04191116161482357
$ws.Cells.Item(4, 26).Value = 0.2414698162729659
$ws.Cells.Item(4, 27).Value = 0.1422328538046669
$ws.Cells.Item(4, 28).Value = 0.07464280583377668
$ws.Cells.Item(4, 29).Value = 0.06759004797089017
$ws.Cells.Item(4, 30).Value = 3.59
$ws.Cells.Item(4, 32).Value = 3.59
$ws.Cells.Item(4, 33).Value = 3.561
$ws.Cells.Item(4, 34).Value = 0.2017987633501968
$ws.Cells.Item(4, 35).Value = 0.2347939829954218
$ws.Cells.Item(4, 36).Value = 0.2004954675975452
$ws.Cells.Item(4, 37).Value = 0.2333398859838805
$ws.Cells.Item(4, 38).Value = 0.246
$ws.Cells.Item(4, 39).Value = 0.246
$ws.Cells.Item(4, 40).Value = 1.879581151832461
$ws.Cells.Item(4, 41).Value = 7.479674796747968
$ws.Cells.Item(4, 42).Value = 1.864397905759162
$ws.Cells.Item(4, 43).Value = 7.479674796747968

# Row 5 writes
$ws.Cells.Item(5, 2).Value = "JS Investments Limited (KASE:JSIL)"
$ws.Cells.Item(5, 4).Value = -0.0599
$ws.Cells.Item(5, 5).Value = -0.197
$ws.Cells.Item(5, 7).Value = -0.140072202166065
$ws.Cells.Item(5, 8).Value = -0.140072202166065
$ws.Cells.Item(5, 9).Value = 0.2227436823104693
$ws.Cells.Item(5, 10).Value = 0.2227436823104693
$ws.Cells.Item(5, 11).Value = 0.818
$ws.Cells.Item(5, 12).Value = 0.2953068592057762
$ws.Cells.Item(5, 13).Value = 1.48
$ws.Cells.Item(5, 14).Value = 0.1649944258639911
$ws.Cells.Item(5, 15).Value = 1.809290953545232
$ws.Cells.Item(5, 16).Value = 1.48
$ws.Cells.Item(5, 17).Value = 0.1649944258639911
$ws.Cells.Item(5, 18).Value = 1.809290953545232
$ws.Cells.Item(5, 21).Value = 0.06
$ws.Cells.Item(5, 22).Value = 0.00668896321070234
$ws.Cells.Item(5, 23).Value = 0.06704918032786886
$ws.Cells.Item(5, 24).Value = 0.07830147531764248
$ws.Cells.Item(5, 25).Value = -0.01125229498977362
$ws.Cells.Item(5, 26).Value = 0.1999566880820039
$ws.Cells.Item(5, 27).Value = 0.04453908900599148
$ws.Cells.Item(5, 28).Value = 0.07404041035876871
$ws.Cells.Item(5, 29).Value = -0.02950132135277723
$ws.Cells.Item(5, 30).Value = 1.92
$ws.Cells.Item(5, 32).Value = 1.92
$ws.Cells.Item(5, 33).Value = 1.86
$ws.Cells.Item(5, 34).Value = 0.1763085399449036
$ws.Cells.Item(5, 35).Value = 0.1521394611727417
$ws.Cells.Item(5, 36).Value = 0.1717451523545706
$ws.Cells.Item(5, 37).Value = 0.1480891719745223
$ws.Cells.Item(5, 38).Value = 0.217
$ws.Cells.Item(5, 39).Value = 0.114
$ws.Cells.Item(5, 40).Value = 2.935779816513761
$ws.Cells.Item(5, 41).Value = 2.84331797235023
$ws.Cells.Item(5, 42).Value = 2.844036697247706
$ws.Cells.Item(5, 43).Value = 5.412280701754386

# Row 6 writes
$ws.Cells.Item(6, 2).Value = "Sindh Modaraba (KASE:SINDM)"
$ws.Cells.Item(6, 11).Value = 0.642
$ws.Cells.Item(6, 12).Value = 0.6557711950970379
$ws.Cells.Item(6, 13).Value = 0.365
$ws.Cells.Item(6, 14).Value = 0.1489795918367347
$ws.Cells.Item(6, 15).Value = 0.5685358255451713
$ws.Cells.Item(6, 16).Value = 0.365
$ws.Cells.Item(6, 17).Value = 0.1489795918367347
$ws.Cells.Item(6, 18).Value = 0.5685358255451713
$ws.Cells.Item(6, 20).Value = 0.0
$ws.Cells.Item(6, 21).Value = 5.58
$ws.Cells.Item(6, 22).Value = 2.277551020408163
$ws.Cells.Item(6, 23).Value = 0.09771689497716894
$ws.Cells.Item(6, 24).Value = 0.06987381421014414
$ws.Cells.Item(6, 25).Value = 0.0278430807670248
$ws.Cells.Item(6, 26).Value = -48.95000000000104
$ws.Cells.Item(6, 27).Value = -0.0
$ws.Cells.Item(6, 28).Value = 0.06987381421014414
$ws.Cells.Item(6, 29).Value = -0.06987381421014414
$ws.Cells.Item(6, 30).Value = 0.0
$ws.Cells.Item(6, 32).Value = 0.0
$ws.Cells.Item(6, 33).Value = -5.58
$ws.Cells.Item(6, 34).Value = 0.0
$ws.Cells.Item(6, 35).Value = 0.0
$ws.Cells.Item(6, 36).Value = 1.782747603833866
$ws.Cells.Item(6, 37).Value = -1.430769230769231
# Row 6 clears
$ws.Cells.Item(6, 4).ClearContents()

# Row 7 writes
$ws.Cells.Item(7, 2).Value = "Unicap Modaraba (KASE:UCAPM)"
$ws.Cells.Item(7, 4).Value = 0.07730000000000001
$ws.Cells.Item(7, 7).Value = 0.0
$ws.Cells.Item(7, 8).Value = 0.0
$ws.Cells.Item(7, 9).Value = 0.0
$ws.Cells.Item(7, 10).Value = 0.0
$ws.Cells.Item(7, 11).Value = -0.019
$ws.Cells.Item(7, 12).Value = -3.8
$ws.Cells.Item(7, 13).Value = -0.0
$ws.Cells.Item(7, 14).Value = -0.0
$ws.Cells.Item(7, 15).Value = 0.0
$ws.Cells.Item(7, 19).Value = 0.0
$ws.Cells.Item(7, 21).Value = 0.0
$ws.Cells.Item(7, 22).Value = 0.0
$ws.Cells.Item(7, 23).Value = -0.08520179372197309
$ws.Cells.Item(7, 24).Value = 0.06987381421014414
$ws.Cells.Item(7, 25).Value = -0.1550756079321172
$ws.Cells.Item(7, 26).Value = 0.02242152466367713
$ws.Cells.Item(7, 27).Value = 0.0
$ws.Cells.Item(7, 28).Value = 0.06987381421014414
$ws.Cells.Item(7, 29).Value = -0.06987381421014414
$ws.Cells.Item(7, 30).Value = 0.0
$ws.Cells.Item(7, 32).Value = 0.0
$ws.Cells.Item(7, 33).Value = 0.0
$ws.Cells.Item(7, 34).Value = 0.0
$ws.Cells.Item(7, 35).Value = 0.0
$ws.Cells.Item(7, 36).Value = 0.0
$ws.Cells.Item(7, 37).Value = 0.0
$ws.Cells.Item(7, 38).Value = 0.0
$ws.Cells.Item(7, 39).Value = 0.0
# Row 7 clears
$ws.Cells.Item(7, 20).ClearContents()
$ws.Cells.Item(7, 40).ClearContents()
$ws.Cells.Item(7, 41).ClearContents()
$ws.Cells.Item(7, 42).ClearContents()
$ws.Cells.Item(7, 43).ClearContents()

# Row 8 writes
$ws.Cells.Item(8, 4).Value = 0.0602
$ws.Cells.Item(8, 5).Value = 0.404
$ws.Cells.Item(8, 7).Value = -0.3223140495867769
$ws.Cells.Item(8, 8).Value = -0.3223140495867769
$ws.Cells.Item(8, 9).Value = -0.1074380165289256
$ws.Cells.Item(8, 10).Value = -0.09208972845336483
$ws.Cells.Item(8, 11).Value = 0.054
$ws.Cells.Item(8, 12).Value = 0.4462809917355372
$ws.Cells.Item(8, 15).Value = -0.0
$ws.Cells.Item(8, 18).Value = -0.0
$ws.Cells.Item(8, 21).Value = 0.209
$ws.Cells.Item(8, 22).Value = 0.19
$ws.Cells.Item(8, 23).Value = 0.0421875
$ws.Cells.Item(8, 24).Value = 0.06987381421014414
$ws.Cells.Item(8, 25).Value = -0.02768631421014415
$ws.Cells.Item(8, 26).Value = 0.1009174311926605
$ws.Cells.Item(8, 27).Value = -0.009293458834743239
$ws.Cells.Item(8, 28).Value = 0.06987381421014414
$ws.Cells.Item(8, 29).Value = -0.07916727304488738
$ws.Cells.Item(8, 33).Value = -0.209
$ws.Cells.Item(8, 36).Value = -0.2345679012345679
$ws.Cells.Item(8, 37).Value = -0.1769686706181202
$ws.Cells.Item(8, 38).Value = 0.001
$ws.Cells.Item(8, 39).Value = -0.007
$ws.Cells.Item(8, 41).Value = -13.0
$ws.Cells.Item(8, 42).Value = 209.0
$ws.Cells.Item(8, 43).Value = 1.857142857142857

# Row 9 writes
$ws.Cells.Item(9, 2).Value = "Arpak International Investments Limited (KASE:ARPAK)"
$ws.Cells.Item(9, 4).Value = 0.0518
$ws.Cells.Item(9, 5).Value = 0.5
$ws.Cells.Item(9, 7).Value = -0.0202020202020202
$ws.Cells.Item(9, 8).Value = -0.0202020202020202
$ws.Cells.Item(9, 9).Value = -5.575757575757576
$ws.Cells.Item(9, 10).Value = -5.484351713859911
$ws.Cells.Item(9, 11).Value = 0.78
$ws.Cells.Item(9, 12).Value = 7.878787878787879
$ws.Cells.Item(9, 13).Value = -0.0
$ws.Cells.Item(9, 14).Value = -0.0
$ws.Cells.Item(9, 15).Value = -0.0
$ws.Cells.Item(9, 16).Value = -0.0
$ws.Cells.Item(9, 17).Value = -0.0
$ws.Cells.Item(9, 18).Value = -0.0
$ws.Cells.Item(9, 21).Value = 0.0
$ws.Cells.Item(9, 22).Value = 0.0
$ws.Cells.Item(9, 23).Value = 0.2015503875968992
$ws.Cells.Item(9, 24).Value = 0.06987381421014414
$ws.Cells.Item(9, 25).Value = 0.1316765733867551
$ws.Cells.Item(9, 26).Value = 0.02559462254395036
$ws.Cells.Item(9, 27).Value = -0.1403699120145117
$ws.Cells.Item(9, 28).Value = 0.06987381421014414
$ws.Cells.Item(9, 29).Value = -0.2102437262246558
$ws.Cells.Item(9, 33).Value = 0.0
$ws.Cells.Item(9, 36).Value = 0.0
$ws.Cells.Item(9, 37).Value = 0.0
$ws.Cells.Item(9, 38).Value = 0.0
$ws.Cells.Item(9, 39).Value = 0.0
$ws.Cells.Item(9, 42).Value = -0.0
# Row 9 clears
$ws.Cells.Item(9, 20).ClearContents()
$ws.Cells.Item(9, 41).ClearContents()
$ws.Cells.Item(9, 43).ClearContents()
